# Insert a new data row at row 96 (pushes existing rows 96:133 down to 97:134)
# and populate it with the new weekly price record, per the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96, 3).Value = "Ñuble"
$ws.Cells.Item(96, 4).Value = 44567
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = 100112028
$ws.Cells.Item(96, 7).Value = "Sandia"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 600
$ws.Cells.Item(96, 11).Value = 1800
$ws.Cells.Item(96, 12).Value = 2200
$ws.Cells.Item(96, 13).Value = 2000
$ws.Cells.Item(96, 14).Value = "`$/unidad"
$ws.Cells.Item(96, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(96, 16).Value = 2000
$ws.Cells.Item(96, 17).Value = 1
$ws.Cells.Item(96, 18).Value = "Hortaliza"
